$wb = $excel.ActiveWorkbook

# --- Sheet "Data": append new weekly observation row ---
$data = $wb.Worksheets.Item("Data")
# Carry the prior row's formatting (date number format, border, bold, centering)
# down onto the newly appended row before filling in its values.
$data.Range("A94:B94").Copy($data.Range("A95:B95"))
$data.Cells.Item(95, 1).Value = 45126
$data.Cells.Item(95, 2).Value = 8274.552

# --- Sheet "SeriesInfo": refresh metadata fields pulled from FRED ---
$info = $wb.Worksheets.Item("SeriesInfo")

# These values (e.g. "2023-07-24") look like dates, so a plain .Value
# assignment would get auto-converted to a date serial + date number
# format. Force text-cell semantics by formatting as Text first, then
# strip the leftover "@" number format back off by re-pasting the
# formatting of an untouched plain-text cell (B2) on top, so the cell
# keeps its original (default) style.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $info.Range("B2").Copy()
    $range.PasteSpecial(-4122)
}

Set-TextValue $info.Range("B3") "2023-07-24"
Set-TextValue $info.Range("B4") "2023-07-24"
Set-TextValue $info.Range("B7") "2023-07-19"
Set-TextValue $info.Range("B14") "2023-07-20 15:33:32-05"
